$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @{"A"="ECs"; "B"="Sema6a"; "C"="Plxna2"; "D"="ECs"; "E"=2; "F"=0.6666666666666666; "G"=42.665376; "H"=127.996128; "I"=0.5598448706913429; "J"=0.5598448706913429; "K"=3; "L"=1; "M"=29.420614; "N"=88.261842; "O"=0.5865186809777162; "P"=0.5865186809777162; "Q"=1255.241558460864; "R"=11297.17402614778; "S"=0.3283594751100266; "T"=0.3283594751100266}
  @{"A"="ECs"; "B"="Sema6a"; "C"="Plxna2"; "D"="FAPs"; "E"=2; "F"=0.6666666666666666; "G"=42.665376; "H"=127.996128; "I"=0.5598448706913429; "J"=0.5598448706913429; "K"=3; "L"=1; "M"=4.080312; "N"=12.240936; "O"=0.08134361887272465; "P"=0.08134361887272466; "Q"=174.088045677312; "R"=1566.792411095808; "S"=0.04553980778936641; "T"=0.04553980778936642}
  @{"A"="ECs"; "B"="Sema6a"; "C"="Plxna2"; "D"="sCs"; "E"=2; "F"=0.6666666666666666; "G"=42.665376; "H"=127.996128; "I"=0.5598448706913429; "J"=0.5598448706913429; "K"=3; "L"=1; "M"=16.660501; "N"=49.981503; "O"=0.3321377001495591; "P"=0.3321377001495591; "Q"=710.826539513376; "R"=6397.438855620384; "S"=0.1859455877919499; "T"=0.1859455877919499}
  @{"A"="FAPs"; "B"="Sema6a"; "C"="Plxna2"; "D"="ECs"; "E"=3; "F"=1; "G"=4.56506; "H"=13.69518; "I"=0.05990162668197795; "J"=0.05990162668197797; "K"=3; "L"=1; "M"=29.420614; "N"=88.261842; "O"=0.5865186809777162; "P"=0.5865186809777162; "Q"=134.30686814684; "R"=1208.76181332156; "S"=0.03513342306993328; "T"=0.03513342306993329}
  @{"A"="FAPs"; "B"="Sema6a"; "C"="Plxna2"; "D"="FAPs"; "E"=3; "F"=1; "G"=4.56506; "H"=13.69518; "I"=0.05990162668197795; "J"=0.05990162668197797; "K"=3; "L"=1; "M"=4.080312; "N"=12.240936; "O"=0.08134361887272465; "P"=0.08134361887272466; "Q"=18.62686909872; "R"=167.64182188848; "S"=0.004872615090675048; "T"=0.00487261509067505}
  @{"A"="FAPs"; "B"="Sema6a"; "C"="Plxna2"; "D"="sCs"; "E"=3; "F"=1; "G"=4.56506; "H"=13.69518; "I"=0.05990162668197795; "J"=0.05990162668197797; "K"=3; "L"=1; "M"=16.660501; "N"=49.981503; "O"=0.3321377001495591; "P"=0.3321377001495591; "Q"=76.05618669505999; "R"=684.50568025554; "S"=0.01989558852136962; "T"=0.01989558852136963}
  @{"A"="sCs"; "B"="Sema6a"; "C"="Plxna2"; "D"="ECs"; "E"=3; "F"=1; "G"=28.97884666666667; "H"=86.93654; "I"=0.3802535026266791; "J"=0.3802535026266792; "K"=3; "L"=1; "M"=29.420614; "N"=88.261842; "O"=0.5865186809777162; "P"=0.5865186809777162; "Q"=852.5754619451867; "R"=7673.17915750668; "S"=0.2230257827977564; "T"=0.2230257827977564}
  @{"A"="sCs"; "B"="Sema6a"; "C"="Plxna2"; "D"="FAPs"; "E"=3; "F"=1; "G"=28.97884666666667; "H"=86.93654; "I"=0.3802535026266791; "J"=0.3802535026266792; "K"=3; "L"=1; "M"=4.080312; "N"=12.240936; "O"=0.08134361887272465; "P"=0.08134361887272466; "Q"=118.24273580016; "R"=1064.18462220144; "S"=0.03093119599268319; "T"=0.0309311959926832}
  @{"A"="sCs"; "B"="Sema6a"; "C"="Plxna2"; "D"="sCs"; "E"=3; "F"=1; "G"=28.97884666666667; "H"=86.93654; "I"=0.3802535026266791; "J"=0.3802535026266792; "K"=3; "L"=1; "M"=16.660501; "N"=49.981503; "O"=0.3321377001495591; "P"=0.3321377001495591; "Q"=482.8021038688466; "R"=4345.21893481962; "S"=0.1262965238362395; "T"=0.1262965238362396}
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$colIndex = @{}
for ($i = 0; $i -lt $columns.Length; $i++) {
  $colIndex[$columns[$i]] = $i + 1
}

$r = 2
foreach ($row in $data) {
  foreach ($c in $columns) {
    $ws.Cells.Item($r, $colIndex[$c]).Value = $row[$c]
  }
  $r++
}

Write-Output "done"